$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Duplicate existing requirement "blocks" to create the new ones,
#    so the new rows inherit the exact same fills/borders/merges
#    already used by their sibling blocks.
# ------------------------------------------------------------------

# RF-003 block (rows 16-21) <- copy of RF-002 block (rows 9-14)
$ws.Range("B9:F14").Copy($ws.Range("B16:F21"))

# RF-004 block (rows 23-28) <- copy of RF-002 block (rows 9-14)
$ws.Range("B9:F14").Copy($ws.Range("B23:F28"))

# RNF-005 block (rows 14-15) <- copy of RNF-004 block (rows 11-12)
$ws.Range("H11:L12").Copy($ws.Range("H14:L15"))

# RI-004 block (rows 11-12) <- copy of RI-003 block (rows 8-9)
$ws.Range("N8:R9").Copy($ws.Range("N11:R12"))

# RI-005 block (rows 14-15) <- copy of RI-003 block (rows 8-9)
$ws.Range("N8:R9").Copy($ws.Range("N14:R15"))

# ------------------------------------------------------------------
# 2) Fill in the real text for the newly created blocks.
# ------------------------------------------------------------------

# RF-003: Mudar linguagem
$ws.Range("B16").Value = "RF-003: Mudar linguagem"
$ws.Range("C17").Value = "Esta funcionalidade permite ao utilizador alterar entre português e inglês o conteúdo apresentado."
$ws.Range("C18").Value = "Urgente - NÃO`nImportante -SIM"
$ws.Range("C19").Value = "Utilizador anónimo"
$ws.Range("C20").Value = "N/A"
$ws.Range("C21").Value = "N/A"

# RF-004: Voltar para o topo
$ws.Range("B23").Value = "RF-004: Voltar para o topo"
$ws.Range("C24").Value = "Esta funcionalidade permite ao utilizador através de um botão voltar para o topo da página."
$ws.Range("C25").Value = "Urgente - NÃO`nImportante -NÃO"
$ws.Range("C26").Value = "Utilizador anónimo"
$ws.Range("C27").Value = "N/A"
$ws.Range("C28").Value = "Apenas deverá ficar visível apos serem movidos no eixo Y 300 pixeis."

# RNF-005: Compatibilidade
$ws.Range("H14").Value = "RNF-005: Compatibilidade"
$ws.Range("I15").Value = "A aplicação deve ser completamente funcional sem muita discrepância de web browser para web browser."

# RI-004:  Cor
$ws.Range("N11").Value = "RI-004:  Cor"
$ws.Range("O12").Value = "As cores devem estar em conformidade da basicamente, ou seja, branco(#ffffff) e azul(#0c5eac)."

# RI-005:  Formato de imagens
$ws.Range("N14").Value = "RI-005:  Formato de imagens"
$ws.Range("O15").Value = "O formato das imagens deve ser do tipo webp."

# ------------------------------------------------------------------
# 3) Edits to pre-existing cells.
# ------------------------------------------------------------------

# RI-001 Tipografia: description now also mentions the typeface
$ws.Range("O3").Value = "O tipo de letra deve ser Poppins e o tamanho da fonte não deve ser menor que 0.875rem/14px."

# RF-001 / RF-002 Autores: "Qualquer utilizador" -> "Utilizador anónimo"
$ws.Range("C5").Value = "Utilizador anónimo"
$ws.Range("C12").Value = "Utilizador anónimo"

# RNF-004 Desempenho: description now also mentions SEO
$ws.Range("I12").Value = "Devem ser tidos em conta formas de melhorar de a performance da aplicação e o SEO, sejá através do uso de imagens com formato .webp e loading do tipo lazy, a bundles de .css, .js, purgados e minificados de forma a reduzir o payload de cada request, etc…"

# ------------------------------------------------------------------
# 4) Misc view state to match the saved workbook.
# ------------------------------------------------------------------
$ws.Range("N21").Select()
